{"js": "// Replace the date line and the 25 multiplication-equation cells with the\n// new values from the diff. Every \"old\" string is unique in the document,\n// so a scoped Find (body.search) + format-preserving Replace\n// (range.insertText(text, \"Replace\")) is safe for each pair.\nconst replacements = [\n  [\"2024-04-26 Friday\", \"2024-04-27 Saturday\"],\n  [\"513\u00d74=2052\", \"704\u00d77=4928\"],\n  [\"899\u00d77=6293\", \"867\u00d72=1734\"],\n  [\"394\u00d79=3546\", \"790\u00d73=2370\"],\n  [\"790\u00d72=1580\", \"678\u00d77=4746\"],\n  [\"707\u00d78=5656\", \"511\u00d72=1022\"],\n  [\"802\u00d75=4010\", \"525\u00d77=3675\"],\n  [\"528\u00d76=3168\", \"862\u00d78=6896\"],\n  [\"438\u00d72=876\", \"838\u00d79=7542\"],\n  [\"978\u00d75=4890\", \"339\u00d79=3051\"],\n  [\"123\u00d74=492\", \"819\u00d73=2457\"],\n  [\"702\u00d79=6318\", \"727\u00d76=4362\"],\n  [\"299\u00d74=1196\", \"505\u00d77=3535\"],\n  [\"200\u00d73=600\", \"144\u00d72=288\"],\n  [\"337\u00d73=1011\", \"655\u00d77=4585\"],\n  [\"651\u00d77=4557\", \"147\u00d79=1323\"],\n  [\"715\u00d78=5720\", \"182\u00d79=1638\"],\n  [\"510\u00d75=2550\", \"733\u00d74=2932\"],\n  [\"517\u00d73=1551\", \"413\u00d79=3717\"],\n  [\"518\u00d76=3108\", \"359\u00d79=3231\"],\n  [\"701\u00d77=4907\", \"894\u00d76=5364\"],\n  [\"539\u00d73=1617\", \"993\u00d76=5958\"],\n  [\"308\u00d77=2156\", \"329\u00d77=2303\"],\n  [\"916\u00d73=2748\", \"283\u00d74=1132\"],\n  [\"519\u00d75=2595\", \"786\u00d72=1572\"],\n  [\"816\u00d74=3264\", \"408\u00d78=3264\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and the 25 multiplication-equation cells with the\n# new values from the diff. Every \"old\" string is unique in the document,\n# so a plain Find/Replace-All over the whole document content is safe for\n# each pair and naturally preserves the existing run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-04-26 Friday\", \"2024-04-27 Saturday\"),\n    @(\"513\u00d74=2052\", \"704\u00d77=4928\"),\n    @(\"899\u00d77=6293\", \"867\u00d72=1734\"),\n    @(\"394\u00d79=3546\", \"790\u00d73=2370\"),\n    @(\"790\u00d72=1580\", \"678\u00d77=4746\"),\n    @(\"707\u00d78=5656\", \"511\u00d72=1022\"),\n    @(\"802\u00d75=4010\", \"525\u00d77=3675\"),\n    @(\"528\u00d76=3168\", \"862\u00d78=6896\"),\n    @(\"438\u00d72=876\", \"838\u00d79=7542\"),\n    @(\"978\u00d75=4890\", \"339\u00d79=3051\"),\n    @(\"123\u00d74=492\", \"819\u00d73=2457\"),\n    @(\"702\u00d79=6318\", \"727\u00d76=4362\"),\n    @(\"299\u00d74=1196\", \"505\u00d77=3535\"),\n    @(\"200\u00d73=600\", \"144\u00d72=288\"),\n    @(\"337\u00d73=1011\", \"655\u00d77=4585\"),\n    @(\"651\u00d77=4557\", \"147\u00d79=1323\"),\n    @(\"715\u00d78=5720\", \"182\u00d79=1638\"),\n    @(\"510\u00d75=2550\", \"733\u00d74=2932\"),\n    @(\"517\u00d73=1551\", \"413\u00d79=3717\"),\n    @(\"518\u00d76=3108\", \"359\u00d79=3231\"),\n    @(\"701\u00d77=4907\", \"894\u00d76=5364\"),\n    @(\"539\u00d73=1617\", \"993\u00d76=5958\"),\n    @(\"308\u00d77=2156\", \"329\u00d77=2303\"),\n    @(\"916\u00d73=2748\", \"283\u00d74=1132\"),\n    @(\"519\u00d75=2595\", \"786\u00d72=1572\"),\n    @(\"816\u00d74=3264\", \"408\u00d78=3264\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n\nWrite-Output \"done\"\n"}
